# Update "Correspond Handoff Datetime" (E4) and "Correspond Handback DateTime" (H4)
# timestamps on the zh-cn and de-de report sheets to reflect the regenerated
# handback status report.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-22 00:40:47"
$wsZh.Range("H4").Value = "2016-03-22 00:41:13"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-22 00:40:53"
$wsDe.Range("H4").Value = "2016-03-22 00:41:26"
